$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price values in column D look numeric (e.g. "1.011") but must remain
# plain text, matching the original inline-string cells. Assigning a
# numeric-looking string via .Value lets Excel auto-convert it to a real
# number, so we briefly force the Text number format for each D cell,
# write the string, then restore the default "Normal" style/format so
# the cell formatting is left exactly as it was originally.

$d = $ws.Range("D2")
$d.NumberFormat = "@"
$d.Value = "29.120.67"
$d.Style = "Normal"
$ws.Range("E2").Value = "  -3.89%  "

$d = $ws.Range("D3")
$d.NumberFormat = "@"
$d.Value = "1.964.83"
$d.Style = "Normal"
$ws.Range("E3").Value = "  -6.41%  "

$ws.Range("E4").Value = "  +0.86%  "

$d = $ws.Range("D5")
$d.NumberFormat = "@"
$d.Value = "327.40"
$d.Style = "Normal"
$ws.Range("E5").Value = "  -4.65%  "

$d = $ws.Range("D6")
$d.NumberFormat = "@"
$d.Value = "1.010"
$d.Style = "Normal"
$ws.Range("E6").Value = "  +0.75%  "

$d = $ws.Range("D7")
$d.NumberFormat = "@"
$d.Value = "0.4995"
$d.Style = "Normal"
$ws.Range("E7").Value = "  -5.13%  "

$d = $ws.Range("D8")
$d.NumberFormat = "@"
$d.Value = "0.4218"
$d.Style = "Normal"
$ws.Range("E8").Value = "  -4.18%  "

$d = $ws.Range("D9")
$d.NumberFormat = "@"
$d.Value = "54.26"
$d.Style = "Normal"
$ws.Range("E9").Value = "  -1.49%  "

$d = $ws.Range("D10")
$d.NumberFormat = "@"
$d.Value = "0.09091"
$d.Style = "Normal"
$ws.Range("E10").Value = "  -3.16%  "

$d = $ws.Range("D11")
$d.NumberFormat = "@"
$d.Value = "1.099"
$d.Style = "Normal"
$ws.Range("E11").Value = "  -6.66%  "

$d = $ws.Range("D12")
$d.NumberFormat = "@"
$d.Value = "23.02"
$d.Style = "Normal"
$ws.Range("E12").Value = "  -7.21%  "

$d = $ws.Range("D13")
$d.NumberFormat = "@"
$d.Value = "1.982.60"
$d.Style = "Normal"
$ws.Range("E13").Value = "  -4.42%  "

$d = $ws.Range("D14")
$d.NumberFormat = "@"
$d.Value = "7.882"
$d.Style = "Normal"
$ws.Range("E14").Value = "  -7.99%  "

$d = $ws.Range("D15")
$d.NumberFormat = "@"
$d.Value = "6.435"
$d.Style = "Normal"
$ws.Range("E15").Value = "  -6.48%  "

$d = $ws.Range("D16")
$d.NumberFormat = "@"
$d.Value = "1.011"
$d.Style = "Normal"
$ws.Range("E16").Value = "  +0.78%  "

$d = $ws.Range("D17")
$d.NumberFormat = "@"
$d.Value = "0.00001101"
$d.Style = "Normal"
$ws.Range("E17").Value = "  -4.91%  "

$d = $ws.Range("D18")
$d.NumberFormat = "@"
$d.Value = "91.26"
$d.Style = "Normal"
$ws.Range("E18").Value = "  -10.04%  "

$d = $ws.Range("D19")
$d.NumberFormat = "@"
$d.Value = "0.06676"
$d.Style = "Normal"
$ws.Range("E19").Value = "  -0.79%  "

$d = $ws.Range("D20")
$d.NumberFormat = "@"
$d.Value = "19.24"
$d.Style = "Normal"
$ws.Range("E20").Value = "  -9.08%  "

$d = $ws.Range("D21")
$d.NumberFormat = "@"
$d.Value = "1.009"
$d.Style = "Normal"
$ws.Range("E21").Value = "  +0.75%  "

$d = $ws.Range("D22")
$d.NumberFormat = "@"
$d.Value = "5.946"
$d.Style = "Normal"
$ws.Range("E22").Value = "  -7.91%  "

$d = $ws.Range("D23")
$d.NumberFormat = "@"
$d.Value = "29.137.06"
$d.Style = "Normal"
$ws.Range("E23").Value = "  -3.88%  "

$ws.Range("E24").Value = "  -3.75%  "

$d = $ws.Range("D25")
$d.NumberFormat = "@"
$d.Value = "2.295"
$d.Style = "Normal"
$ws.Range("E25").Value = "  -1.19%  "

$d = $ws.Range("D26")
$d.NumberFormat = "@"
$d.Value = "157.04"
$d.Style = "Normal"
$ws.Range("E26").Value = "  -3.30%  "

$d = $ws.Range("D27")
$d.NumberFormat = "@"
$d.Value = "20.64"
$d.Style = "Normal"
$ws.Range("E27").Value = "  -5.45%  "

$d = $ws.Range("D28")
$d.NumberFormat = "@"
$d.Value = "6.213"
$d.Style = "Normal"
$ws.Range("E28").Value = "  -10.79%  "

$ws.Range("E29").Value = "  -10.28%  "

$d = $ws.Range("D30")
$d.NumberFormat = "@"
$d.Value = "127.27"
$d.Style = "Normal"
$ws.Range("E30").Value = "  -4.76%  "

$d = $ws.Range("D31")
$d.NumberFormat = "@"
$d.Value = "1.041"
$d.Style = "Normal"
$ws.Range("E31").Value = "  -8.83%  "

$d = $ws.Range("D32")
$d.NumberFormat = "@"
$d.Value = "0.09855"
$d.Style = "Normal"
$ws.Range("E32").Value = "  -6.52%  "

$d = $ws.Range("D33")
$d.NumberFormat = "@"
$d.Value = "1.528"
$d.Style = "Normal"
$ws.Range("E33").Value = "  -8.83%  "

$d = $ws.Range("D34")
$d.NumberFormat = "@"
$d.Value = "5.823"
$d.Style = "Normal"
$ws.Range("E34").Value = "  -7.17%  "

$d = $ws.Range("D35")
$d.NumberFormat = "@"
$d.Value = "3.711"
$d.Style = "Normal"
$ws.Range("E35").Value = "  -4.30%  "

$d = $ws.Range("D36")
$d.NumberFormat = "@"
$d.Value = "0.02429"
$d.Style = "Normal"
$ws.Range("E36").Value = "  -8.09%  "

$d = $ws.Range("D37")
$d.NumberFormat = "@"
$d.Value = "9.008"
$d.Style = "Normal"
$ws.Range("E37").Value = "  -10.83%  "

$d = $ws.Range("D38")
$d.NumberFormat = "@"
$d.Value = "0.06346"
$d.Style = "Normal"
$ws.Range("E38").Value = "  -6.41%  "

$d = $ws.Range("D39")
$d.NumberFormat = "@"
$d.Value = "1.288"
$d.Style = "Normal"
$ws.Range("E39").Value = "  -4.32%  "

$d = $ws.Range("D40")
$d.NumberFormat = "@"
$d.Value = "0.6443"
$d.Style = "Normal"
$ws.Range("E40").Value = "  -7.58%  "

$d = $ws.Range("D41")
$d.NumberFormat = "@"
$d.Value = "11.47"
$d.Style = "Normal"
$ws.Range("E41").Value = "  -9.60%  "

$d = $ws.Range("D42")
$d.NumberFormat = "@"
$d.Value = "0.2006"
$d.Style = "Normal"
$ws.Range("E42").Value = "  -9.55%  "

$d = $ws.Range("D43")
$d.NumberFormat = "@"
$d.Value = "1.009"
$d.Style = "Normal"
$ws.Range("E43").Value = "  +0.69%  "

$d = $ws.Range("D44")
$d.NumberFormat = "@"
$d.Value = "0.6221"
$d.Style = "Normal"
$ws.Range("E44").Value = "  -8.31%  "

$d = $ws.Range("D45")
$d.NumberFormat = "@"
$d.Value = "13.40"
$d.Style = "Normal"
$ws.Range("E45").Value = "  -7.46%  "

$d = $ws.Range("D46")
$d.NumberFormat = "@"
$d.Value = "2.179"
$d.Style = "Normal"
$ws.Range("E46").Value = "  -6.44%  "

$d = $ws.Range("D47")
$d.NumberFormat = "@"
$d.Value = "1.296"
$d.Style = "Normal"
$ws.Range("E47").Value = "  -1.06%  "

$d = $ws.Range("D48")
$d.NumberFormat = "@"
$d.Value = "3.469"
$d.Style = "Normal"
$ws.Range("E48").Value = "  -4.61%  "

$d = $ws.Range("D49")
$d.NumberFormat = "@"
$d.Value = "0.00000000332"
$d.Style = "Normal"
$ws.Range("E49").Value = "  -3.32%  "

$d = $ws.Range("D50")
$d.NumberFormat = "@"
$d.Value = "0.06870"
$d.Style = "Normal"
$ws.Range("E50").Value = "  -5.94%  "

$d = $ws.Range("D51")
$d.NumberFormat = "@"
$d.Value = "1.107"
$d.Style = "Normal"
$ws.Range("E51").Value = "  -8.57%  "
